$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1095.3334
$ws.Range("I127").Value = 331.2857
$ws.Range("J127").Value = 2165
$ws.Range("K127").Value = 993.8571000000001
$ws.Range("L127").Value = 6495
$ws.Range("M127").Value = 3966.1429
$ws.Range("N127").Value = -16415
$ws.Range("H129").Value = 1078.1772
$ws.Range("I129").Value = 433.4
$ws.Range("K129").Value = 1300.2
$ws.Range("M129").Value = 3699.8
$ws.Range("H137").Value = 2293.7046
$ws.Range("I137").Value = 1584.2858
$ws.Range("K137").Value = 4752.857400000001
$ws.Range("M137").Value = -2202.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7861.421
$ws.Range("I122").Value = 8179.933
$ws.Range("K122").Value = 24539.799
$ws.Range("M122").Value = -22089.799
$ws.Range("H123").Value = 25426.666
$ws.Range("J123").Value = 25426.666
$ws.Range("L123").Value = 25426.666
$ws.Range("N123").Value = -35226.666
$ws.Range("H124").Value = 35261.8
$ws.Range("J124").Value = 35261.8
$ws.Range("L124").Value = 35261.8
$ws.Range("N124").Value = -45081.8
$ws.Range("H131").Value = 44285.4
$ws.Range("J131").Value = 44285.4
$ws.Range("L131").Value = 44285.4
$ws.Range("N131").Value = -54365.4
$ws.Range("H135").Value = 23087.445
$ws.Range("J135").Value = 23087.445
$ws.Range("L135").Value = 23087.445
$ws.Range("N135").Value = -33227.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6040
$ws.Range("H31").Value = 1754.0494
$ws.Range("I31").Value = 2283.2856
$ws.Range("J31").Value = 1351.3695
$ws.Range("K31").Value = 2283.2856
$ws.Range("L31").Value = 1351.3695
$ws.Range("M31").Value = -1988.2856
$ws.Range("N31").Value = -1941.3695
$ws.Range("H34").Value = 1754.0494
$ws.Range("I34").Value = 2283.2856
$ws.Range("J34").Value = 1351.3695
$ws.Range("K34").Value = 2283.2856
$ws.Range("L34").Value = 1351.3695
$ws.Range("M34").Value = -2081.2856
$ws.Range("N34").Value = -1755.3695
$ws.Range("H132").Value = 295423.1
$ws.Range("I132").Value = 398473.22
$ws.Range("J132").Value = 3447.8333
$ws.Range("K132").Value = 1195419.66
$ws.Range("L132").Value = 10343.4999
$ws.Range("M132").Value = -1192889.66
$ws.Range("N132").Value = -15403.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 696.6667
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 696.6667
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2090.0001
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2314.0001
$ws.Range("H68").Value = 878.65
$ws.Range("I68").Value = 692.18335
$ws.Range("J68").Value = 1158.35
$ws.Range("K68").Value = 2076.55005
$ws.Range("L68").Value = 3475.05
$ws.Range("M68").Value = -1265.55005
$ws.Range("N68").Value = -5097.049999999999
$ws.Range("H71").Value = 878.65
$ws.Range("I71").Value = 692.18335
$ws.Range("J71").Value = 1158.35
$ws.Range("K71").Value = 6229.65015
$ws.Range("L71").Value = 10425.15
$ws.Range("M71").Value = -2173.65015
$ws.Range("N71").Value = -18537.15
$ws.Range("H107").Value = 1187.5211
$ws.Range("I107").Value = 1188.0526
$ws.Range("J107").Value = 1186.909
$ws.Range("K107").Value = 3564.1578
$ws.Range("L107").Value = 3560.727
$ws.Range("M107").Value = -1644.1578
$ws.Range("N107").Value = -7400.727000000001
$ws.Range("H118").Value = 2619.5908
$ws.Range("I118").Value = 628.75
$ws.Range("J118").Value = 3062
$ws.Range("K118").Value = 1886.25
$ws.Range("L118").Value = 9186
$ws.Range("M118").Value = -643.25
$ws.Range("N118").Value = -11672
$ws.Range("H129").Value = 2941683.5
$ws.Range("I129").Value = 256.36365
$ws.Range("J129").Value = 8334299.5
$ws.Range("K129").Value = 769.09095
$ws.Range("L129").Value = 25002898.5
$ws.Range("M129").Value = 4230.90905
$ws.Range("N129").Value = -25012898.5
$ws.Range("H137").Value = 23813984
$ws.Range("I137").Value = 4453.6665
$ws.Range("J137").Value = 30307494
$ws.Range("K137").Value = 13360.9995
$ws.Range("L137").Value = 90922482
$ws.Range("M137").Value = -8260.999500000002
$ws.Range("N137").Value = -90932682
$ws.Range("H138").Value = 2180.682
$ws.Range("I138").Value = 760.75
$ws.Range("J138").Value = 3884.6
$ws.Range("K138").Value = 2282.25
$ws.Range("L138").Value = 11653.8
$ws.Range("M138").Value = 2857.75
$ws.Range("N138").Value = -21933.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10134.667
$ws.Range("I5").Value = 10134.667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 10134.667
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -10022.667
$ws.Range("N5").ClearContents()
$ws.Range("H109").Value = 13485
$ws.Range("J109").Value = 13485
$ws.Range("L109").Value = 13485
$ws.Range("N109").Value = -15565
$ws.Range("H126").Value = 3076.111
$ws.Range("I126").Value = 1825
$ws.Range("J126").Value = 4077
$ws.Range("K126").Value = 5475
$ws.Range("L126").Value = 12231
$ws.Range("M126").Value = -3005
$ws.Range("N126").Value = -17171
$ws.Range("H132").Value = 3403.7778
$ws.Range("I132").Value = 2479.6667
$ws.Range("J132").Value = 5252
$ws.Range("K132").Value = 7439.000100000001
$ws.Range("L132").Value = 15756
$ws.Range("M132").Value = -4909.000100000001
$ws.Range("N132").Value = -20816

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3326.7058
$ws.Range("I40").Value = 2805.4443
$ws.Range("J40").Value = 3913.125
$ws.Range("K40").Value = 2805.4443
$ws.Range("L40").Value = 3913.125
$ws.Range("M40").Value = -2669.4443
$ws.Range("N40").Value = -4185.125
$ws.Range("H82").Value = 3520.6
$ws.Range("I82").Value = 800
$ws.Range("J82").Value = 5334.3335
$ws.Range("K82").Value = 800
$ws.Range("L82").Value = 5334.3335
$ws.Range("M82").Value = -439
$ws.Range("N82").Value = -6056.3335
$ws.Range("H85").Value = 3520.6
$ws.Range("I85").Value = 800
$ws.Range("J85").Value = 5334.3335
$ws.Range("K85").Value = 800
$ws.Range("L85").Value = 5334.3335
$ws.Range("M85").Value = 448
$ws.Range("N85").Value = -7830.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 17668.666
$ws.Range("I11").Value = 26500.5
$ws.Range("K11").Value = 26500.5
$ws.Range("M11").Value = -26358.5
$ws.Range("H123").Value = 23530.383
$ws.Range("J123").Value = 23530.383
$ws.Range("L123").Value = 23530.383
$ws.Range("N123").Value = -33330.383
$ws.Range("H136").Value = 2202.3635
$ws.Range("I136").Value = 2224.9167
$ws.Range("J136").Value = 2175.3
$ws.Range("K136").Value = 6674.750100000001
$ws.Range("L136").Value = 6525.900000000001
$ws.Range("M136").Value = -4124.750100000001
$ws.Range("N136").Value = -11625.9
